$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the greeting text in E8 from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Update the selection / active cell to E8 to match the saved view state
$ws.Range("E8").Select()
